$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-26 Thursday" "2025-06-27 Friday"

Replace-Text "18×77=" "64×86="
Replace-Text "31×79=" "21×59="
Replace-Text "94×50=" "97×97="
Replace-Text "21×22=" "63×64="
Replace-Text "96×98=" "84×47="

Replace-Text "59×18=" "52×62="
Replace-Text "25×29=" "67×76="
Replace-Text "93×39=" "59×88="
Replace-Text "40×46=" "34×54="
Replace-Text "59×67=" "19×29="

Replace-Text "83×99=" "78×26="
Replace-Text "82×52=" "42×53="
Replace-Text "20×13=" "83×94="
Replace-Text "67×19=" "78×26="
Replace-Text "38×96=" "49×95="

Replace-Text "79×88=" "53×18="
Replace-Text "23×62=" "78×47="
Replace-Text "94×70=" "11×16="
Replace-Text "58×23=" "25×99="
Replace-Text "62×12=" "45×89="

Replace-Text "53×49=" "36×33="
Replace-Text "26×58=" "60×42="
Replace-Text "38×45=" "20×78="
Replace-Text "70×53=" "98×50="
Replace-Text "13×12=" "86×13="
